$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the existing data row, shifting the
# original row 1 (String/1/1.1/TRUE()/FALSE()/date) down to row 2.
$ws.Rows("1:1").Insert()

# Give the now-empty row 1 cells the same per-column formatting the data
# row uses (General / General / General / "BOOL"e"AN" / "BOOL"e"AN" / d-m-yy)
# so they reuse the existing style slots instead of staying fully blank.
$ws.Range("A1:C1").NumberFormat = "General"
$ws.Range("D1:E1").NumberFormat = '"BOOL"e"AN"'
$ws.Range("F1").NumberFormat = "d/m/yy"

# The TRUE()/FALSE() formulas in the (now shifted) boolean cells become
# literal boolean values formatted with a new Russian True/False display
# format, instead of remaining live formulas.
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = $false
$ws.Range("D2:E2").NumberFormat = '"ИСТИНА";"ИСТИНА";"ЛОЖЬ"'

# Cursor ends up parked on the relocated date cell, F2.
$ws.Range("F2").Select() | Out-Null
